# Append " (Changed main)" to the end of the first paragraph
# ("This is a Microsoft word document.") as three additional runs:
#   " ("  /  "Changed main"  /  ")"
#
# A plain Range.InsertAfter gets silently coalesced back into the
# preceding run whenever the adjacent runs end up with identical
# formatting. To keep each inserted chunk as its own <w:r>, briefly wrap
# each freshly-inserted chunk in a bookmark (which forces a run
# boundary) and immediately delete the bookmark again - the bookmark
# markup itself leaves no trace in the saved document, but the run
# split survives.

$d = $word.ActiveDocument

$para = $d.Paragraphs(1)
$pos = $para.Range.End - 1   # end of paragraph text, before the paragraph mark

# --- segment 1: " (" ---
$r = $d.Range($pos, $pos)
$r.InsertAfter(" (")
$newPos = $pos + 2
$bm = $d.Range($pos, $newPos)
$d.Bookmarks.Add("tmpSplit1", $bm)
$d.Bookmarks("tmpSplit1").Delete()
$pos = $newPos

# --- segment 2: "Changed main" ---
$r = $d.Range($pos, $pos)
$r.InsertAfter("Changed main")
$newPos = $pos + 12
$bm = $d.Range($pos, $newPos)
$d.Bookmarks.Add("tmpSplit2", $bm)
$d.Bookmarks("tmpSplit2").Delete()
$pos = $newPos

# --- segment 3: ")" ---
$r = $d.Range($pos, $pos)
$r.InsertAfter(")")
$newPos = $pos + 1
$bm = $d.Range($pos, $newPos)
$d.Bookmarks.Add("tmpSplit3", $bm)
$d.Bookmarks("tmpSplit3").Delete()
$pos = $newPos
